$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for the data block B3:G6 based on the diff.
# Rows/cols that are unchanged keep their original value; this just
# rewrites the full block to the new, correct values.
$values = @{
    "B3" = 3; "C3" = 2; "D3" = 4; "E3" = 3; "F3" = 4; "G3" = 4
    "B4" = 3; "C4" = 3; "D4" = 2; "E4" = 3; "F4" = 2; "G4" = 3
    "B5" = 4; "C5" = 2; "D5" = 4; "E5" = 2; "F5" = 4; "G5" = 3
    "B6" = 3; "C6" = 2; "D6" = 2; "E6" = 4; "F6" = 3; "G6" = 4
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
